# "add array data type"
#
# This lesson workbook gained a note about the Array data type:
#   - cell N81 (next to the "OBJECT METHOD THIS" section header row) gets a
#     stray single-space value — an incidental edit left over from the author
#     typing a note there.
#   - cell A193, in the "Vòng lặp" (loops) section, had a typo fixed:
#     "Các cũ ..." -> "Cách cũ ..." ("Các cũ" was a typo for "Cách cũ",
#     meaning "the old way").
#
# Everything else in the underlying OOXML (shared-string index shuffling,
# dimension growing to column N, mergeCells re-ordering, etc.) is a
# mechanical side effect of Excel rewriting the shared-strings table and
# sheet XML on save, so it falls out naturally once the two content edits
# below are made through the Excel object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New cell N81: a single space, styled like the surrounding text
#        (Times New Roman 12pt, same as A81:A86 in this block). ---
$ws.Range("N81").Value = " "
$ws.Range("N81").Font.Name = "Times New Roman"
$ws.Range("N81").Font.Size = 12

# --- 2. Fix the typo in A193: "Các cũ" -> "Cách cũ" ---
# The source text uses non-breaking spaces (U+00A0) between words (matching
# the rest of this workbook's indentation style), so rebuild it the same way
# rather than with plain ASCII spaces.
$nbsp = [char]0x00A0
$newText = "$nbsp$nbsp$nbsp$nbsp" + `
    "Cách${nbsp}cũ${nbsp}chúng${nbsp}ta${nbsp}dùng${nbsp}vòng${nbsp}lặp${nbsp}for(),${nbsp}cách${nbsp}mới${nbsp}dùng${nbsp}for(of);"
$ws.Range("A193").Value = $newText

# --- 3. Cosmetic: move the active selection to where the author ended up. ---
$ws.Range("J195").Select()
